$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; A="SKU_1"; B="DC"; C=6; D=0.98; E=21.11; F=15.17; G=320859232292.6; H=-0.1; I=49.49; J=1641; K=93.58; Q="Rule-based SS"},
    @{Row=3; A="SKU_1"; B="Store"; C=6; D=0.98; E=21.96; F=14.93; G=289695307847.15; H=6.97; I=58.44; J=1641; K=78.63; Q="Rule-based SS"},
    @{Row=4; A="SKU_1"; B="WH"; C=9; D=0.98; E=27.86; F=18.79; G=263162705836.89; H=12.28; I=73.24; J=1641; K=104.53; Q="Rule-based SS"},
    @{Row=5; A="SKU_2"; B="DC"; C=19; D=0.98; E=20.65; F=12.75; G=106020719154.47; H=5.88; I=45.95; J=1641; K=152.29; Q="Rule-based SS"},
    @{Row=6; A="SKU_2"; B="Store"; C=12; D=0.9; E=27.37; F=17.94; G=686142596112.72; H=12.41; I=73.11; J=1641; K=75.98; Q="Rule-based SS"},
    @{Row=7; A="SKU_2"; B="WH"; C=10; D=0.98; E=21.54; F=13.49; G=137093235897.45; H=-3.61; I=45.29; J=1641; K=124.81; Q="Rule-based SS"},
    @{Row=8; A="SKU_3"; B="DC"; C=10; D=0.98; E=18.25; F=12.84; G=340694698452.79; H=-2.2; I=44.78; J=1641; K=121.43; Q="Rule-based SS"},
    @{Row=9; A="SKU_3"; B="Store"; C=8; D=0.95; E=25.25; F=14.13; G=17068860496.06; H=-1.53; I=38.31; J=1641; K=103.58; Q="Rule-based SS"},
    @{Row=10; A="SKU_3"; B="WH"; C=17; D=0.9; E=22.15; F=14.02; G=42078001262.2; H=-0.49; I=37.34; J=1641; K=110.53; Q="Rule-based SS"},
    @{Row=11; A="SKU_4"; B="DC"; C=18; D=0.95; E=23.44; F=15.79; G=45795246869.18; H=0.61; I=49.98; J=1641; K=125.47; Q="Rule-based SS"},
    @{Row=12; A="SKU_4"; B="Store"; C=18; D=0.9; E=23.56; F=15.38; G=167154174354.58; H=7.91; I=52.57; J=1641; K=92.89; Q="Rule-based SS"},
    @{Row=13; A="SKU_4"; B="WH"; C=8; D=0.95; E=18.88; F=11.87; G=128921389467.27; H=1.71; I=43.4; J=1641; K=78.41; Q="Rule-based SS"},
    @{Row=14; A="SKU_5"; B="DC"; C=8; D=0.98; E=22.05; F=16.53; G=63924436405.61; H=5.85; I=57.42; J=1641; K=105.85; Q="Rule-based SS"},
    @{Row=15; A="SKU_5"; B="Store"; C=19; D=0.98; E=27.52; F=17.69; G=90633759964.73; H=-11.17; I=44.08; J=1641; K=225.56; Q="Rule-based SS"},
    @{Row=16; A="SKU_5"; B="WH"; C=14; D=0.9; E=26.34; F=19.91; G=702199878302.32; H=17.09; I=93.89; J=1641; K=62.49; Q="Rule-based SS"}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
}

# Row 4 previously had hybrid_with_no_var_ss / hybrid_with_var_ss values; clear them now
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 15).Value = ""
